# Q3 Update - 2025
# Refresh the Eswatini (SWA) country-of-asylum data block (rows 439-455):
#  - short-url changes from oKAY71 to Qjm681 for every existing data row
#  - refugee/asylum-seeker figures are refreshed for 2024
#  - two new countries of origin (Malawi, Mozambique) are inserted
#  - a new "Unknown" origin row is appended at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the short-url value (column B) for every existing data row.
# ---------------------------------------------------------------------
$ws.Range("B2:B452").Replace("oKAY71", "Qjm681") | Out-Null

# ---------------------------------------------------------------------
# 2. Extend formatting (styles) of the last existing row down to the
#    3 new rows that will be appended (453-455), matching row 452.
# ---------------------------------------------------------------------
$ws.Range("A452:V452").Copy()
$ws.Range("A453:V455").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Rewrite the full data block for rows 439-455 with the refreshed
#    2024 figures (this also covers the newly inserted Malawi,
#    Mozambique and Unknown rows). Every field in this sheet is stored
#    as text, so the values are written through a scratch range
#    formatted as Text and copied across with values-only paste - this
#    keeps numeric-looking entries (e.g. "201") as text instead of
#    silently becoming numbers.
# ---------------------------------------------------------------------
$rows = @(
  ,@("1","Qjm681","1","438","2024","16","Burundi","BDI","BDI","182","Eswatini","SWA","SWZ","201","231","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","439","2024","39","Cameroon","CMR","CMR","182","Eswatini","SWA","SWZ","5","36","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","440","2024","40","Congo","COB","COG","182","Eswatini","SWA","SWZ","16","27","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","441","2024","41","Dem. Rep. of the Congo","COD","COD","182","Eswatini","SWA","SWZ","569","442","11","0","0","0","0","-","0")
  ,@("1","Qjm681","1","442","2024","56","Eritrea","ERT","ERI","182","Eswatini","SWA","SWZ","5","176","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","443","2024","58","Ethiopia","ETH","ETH","182","Eswatini","SWA","SWZ","52","505","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","444","2024","128","Malawi","MLW","MWI","182","Eswatini","SWA","SWZ","0","5","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","445","2024","131","Mozambique","MOZ","MOZ","182","Eswatini","SWA","SWZ","0","20","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","446","2024","141","Nigeria","NIG","NGA","182","Eswatini","SWA","SWZ","0","16","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","447","2024","161","Rwanda","RWA","RWA","182","Eswatini","SWA","SWZ","205","161","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","448","2024","172","Somalia","SOM","SOM","182","Eswatini","SWA","SWZ","215","733","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","449","2024","182","Eswatini","SWA","SWZ","182","Eswatini","SWA","SWZ","0","0","0","0","0","0","12","-","0")
  ,@("1","Qjm681","1","450","2024","185","Syrian Arab Rep.","SYR","SYR","182","Eswatini","SWA","SWZ","0","5","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","451","2024","186","United Rep. of Tanzania","TAN","TZA","182","Eswatini","SWA","SWZ","0","5","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","452","2024","199","Uganda","UGA","UGA","182","Eswatini","SWA","SWZ","7","6","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","453","2024","262","Unknown ","UKN","UNK","182","Eswatini","SWA","SWZ","0","786","0","0","0","0","0","-","0")
  ,@("1","Qjm681","1","454","2024","214","Zimbabwe","ZIM","ZWE","182","Eswatini","SWA","SWZ","9","5","0","0","0","0","0","-","0")
)

$nrows = $rows.Count
$ncols = $rows[0].Count
$arr = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
  for ($j = 0; $j -lt $ncols; $j++) {
    $arr[$i,$j] = $rows[$i][$j]
  }
}

$firstRow = 439
$lastRow = $firstRow + $nrows - 1
$scratch = $ws.Range("A1000:V" + (1000 + $nrows - 1))
$scratch.NumberFormat = "@"
$scratch.Value = $arr
$scratch.Copy()
$ws.Range("A" + $firstRow + ":V" + $lastRow).PasteSpecial(-4163) | Out-Null  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Update the sheet's used-range dimension now extends to row 455.
# ---------------------------------------------------------------------
Write-Host "Rows written:" $nrows
Write-Host "B2:" $ws.Range("B2").Value2
Write-Host "N439:" $ws.Range("N439").Value2
Write-Host "F445:" $ws.Range("F445").Value2
Write-Host "G445:" $ws.Range("G445").Value2
Write-Host "G454:" $ws.Range("G454").Value2
Write-Host "O454:" $ws.Range("O454").Value2
Write-Host "V455:" $ws.Range("V455").Value2
